# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "24.544.92"
Set-TextValue "E2" "  +3.34%  "
Set-TextValue "D3" "1.692.48"
Set-TextValue "E3" "  +1.65%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "D5" "315.30"
Set-TextValue "E5" "  +1.91%  "
Set-TextValue "E6" "  +0.26%  "
Set-TextValue "D7" "0.3931"
Set-TextValue "E7" "  +1.46%  "
Set-TextValue "D8" "0.3999"
Set-TextValue "E8" "  +1.44%  "
Set-TextValue "D9" "1.519"
Set-TextValue "E9" "  +4.34%  "
Set-TextValue "D10" "1.001"
Set-TextValue "E10" "  +0.20%  "
Set-TextValue "D11" "53.09"
Set-TextValue "E11" "  +5.21%  "
Set-TextValue "D12" "0.08723"
Set-TextValue "E12" "  +0.92%  "
Set-TextValue "D13" "7.169"
Set-TextValue "E13" "  +6.61%  "
Set-TextValue "D14" "23.06"
Set-TextValue "E14" "  +2.15%  "
Set-TextValue "E15" "  +0.08%  "
Set-TextValue "D16" "7.572"
Set-TextValue "E16" "  +4.01%  "
Set-TextValue "D17" "1.692.38"
Set-TextValue "E17" "  +1.63%  "
Set-TextValue "D18" "99.69"
Set-TextValue "E18" "  +0.07%  "
Set-TextValue "D19" "0.07042"
Set-TextValue "E19" "  +3.87%  "
Set-TextValue "D20" "19.56"
Set-TextValue "E20" "  +2.61%  "
Set-TextValue "D21" "6.824"
Set-TextValue "E21" "  +3.11%  "
Set-TextValue "E22" "  +0.17%  "
Set-TextValue "D23" "14.00"
Set-TextValue "E23" "  +1.21%  "
Set-TextValue "D24" "24.530.25"
Set-TextValue "E24" "  +3.31%  "
Set-TextValue "D25" "2.991"
Set-TextValue "E25" "  +6.84%  "
Set-TextValue "D26" "2.317"
Set-TextValue "E26" "  +0.16%  "
Set-TextValue "D27" "22.27"
Set-TextValue "E27" "  +2.69%  "
Set-TextValue "D28" "160.60"
Set-TextValue "E28" "  +0.74%  "
Set-TextValue "D29" "5.212"
Set-TextValue "E29" "  +0.27%  "
Set-TextValue "D30" "133.86"
Set-TextValue "E30" "  +3.40%  "
Set-TextValue "D31" "7.469"
Set-TextValue "E31" "  +12.10%  "
Set-TextValue "D32" "1.879.19"
Set-TextValue "E32" "  +1.75%  "
Set-TextValue "D33" "1.088"
Set-TextValue "E33" "  -2.88%  "
Set-TextValue "B34" "InternetComputer(DFINITY)"
Set-TextValue "C34" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "7.287"
Set-TextValue "E34" "  +10.91%  "
Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.08498"
Set-TextValue "E35" "  -0.09%  "
Set-TextValue "D36" "11.34"
Set-TextValue "E36" "  +8.81%  "
Set-TextValue "D37" "1.955"
Set-TextValue "E37" "  -1.01%  "
Set-TextValue "D38" "0.2703"
Set-TextValue "E38" "  +1.71%  "
Set-TextValue "D39" "14.33"
Set-TextValue "E39" "  -0.76%  "
Set-TextValue "D40" "0.02747"
Set-TextValue "E40" "  +9.45%  "
Set-TextValue "D41" "0.09012"
Set-TextValue "E41" "  +2.64%  "
Set-TextValue "D42" "1.473"
Set-TextValue "E42" "  +1.16%  "
Set-TextValue "D43" "0.7610"
Set-TextValue "E43" "  +1.09%  "
Set-TextValue "D44" "0.7151"
Set-TextValue "E44" "  +2.10%  "
Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "15.31"
Set-TextValue "E45" "  +3.03%  "
Set-TextValue "B46" "NEARProtocol"
Set-TextValue "C46" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D46" "2.514"
Set-TextValue "E46" "  +4.42%  "
Set-TextValue "D47" "4.201"
Set-TextValue "E47" "  +2.44%  "
Set-TextValue "E48" "  +0.28%  "
Set-TextValue "D49" "140.54"
Set-TextValue "E49" "  +1.40%  "
Set-TextValue "D50" "1.304"
Set-TextValue "E50" "  +5.53%  "
Set-TextValue "D51" "0.07987"
Set-TextValue "E51" "  +2.74%  "
